# Update crypto price/volume table (and swap the EnergySwap/Elrond rows 47-48)
# per the "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.517.15'
$ws.Range("E2").Value = '  +0.38%  '
$ws.Range("D3").Value = '1.852.93'
$ws.Range("E3").Value = '  -0.40%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '233.82'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.35%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4707'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.72%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("E9").Value = '  -1.56%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '17.55'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +7.28%  '
$ws.Range("D11").Value = '1.845.66'
$ws.Range("E11").Value = '  -0.82%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07427'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.13%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.054'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.87%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '84.70'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.11%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6274'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.39%  '
$ws.Range("D16").Value = '30.499.31'
$ws.Range("E16").Value = '  +0.52%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '244.36'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +5.60%  '
$ws.Range("E18").Value = '  -0.01%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.70'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.85%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007345'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.97%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.001'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.04%  '
$ws.Range("E22").Value = '  -1.14%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.000'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.26%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.276'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.16%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '162.46'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.83%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '18.08'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.53%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.884'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.84%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.1009'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.30%  '
$ws.Range("E29").Value = '  -1.58%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.038'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.80%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.858'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.18%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.04901'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.31%  '
$ws.Range("E33").Value = '  -1.56%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7052'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.23%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.705'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.37%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.01909'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.57%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.685'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.19%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.8734'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.24%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.977'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.09%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '105.22'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.65%  '
$ws.Range("E41").Value = '  +0.05%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4078'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.10%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.503'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.33%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.255'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.43%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '62.98'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.88%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1200'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.16%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.575'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.25%  '
$ws.Range("B48").Value = 'Elrond'
$ws.Range("C48").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '33.40'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.18%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05533'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.96%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.370'
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3688'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.83%  '
